$wb = $excel.ActiveWorkbook

# Insert a new "is_targeted list" worksheet before the
# "bulk_rna_yield...ssue_unit list" sheet (i.e. as the new sheet #6),
# which pushes the remaining list sheets down by one position.
$beforeSheet = $wb.Worksheets.Item("bulk_rna_yield...ssue_unit list")
$targetList = $wb.Worksheets.Add($beforeSheet)
$targetList.Name = "is_targeted list"

# Populate it with the allowed boolean values as text, matching the
# style used by the other "* list" helper sheets (plain text cells).
$targetList.Range("A1").Value = "'TRUE"
$targetList.Range("A2").Value = "'FALSE"

# Point the is_targeted column's validation at the new list sheet
# instead of the old inline "TRUE,FALSE" formula, updating the
# error title/message to match the list-based validations used
# elsewhere in the workbook.
$ws = $wb.Worksheets.Item("Export as TSV")
$col = $ws.Range("N2:N1048576")
$col.Validation.Formula1 = "='is_targeted list'!`$A`$1:`$A`$2"
$col.Validation.ErrorTitle = "Value must come from list"
$col.Validation.ErrorMessage = "Value must be one of: TRUE / FALSE."
